# Updates Price (D) and Volume(1h) (E) columns per the latest cryptos refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.540.40'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '1.565.99'
$ws.Range("E3").Value = '  -1.58%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Formula = '="211.80"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -1.39%  '
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Formula = '="46.30"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  +5.50%  '
$ws.Range("D9").Formula = '="24.14"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  +0.91%  '
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("E11").Value = '  -1.68%  '
$ws.Range("E12").Value = '  -0.54%  '
$ws.Range("D13").Value = '1.789.75'
$ws.Range("E13").Value = '  -1.55%  '
$ws.Range("D14").Value = '1.568.21'
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("E15").Value = '  -2.18%  '
$ws.Range("D16").Value = '28.536.68'
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("E17").Value = '  -2.91%  '
$ws.Range("D18").Formula = '="61.94"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  -3.30%  '
$ws.Range("D19").Formula = '="227.33"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  -2.33%  '
$ws.Range("D20").Value = '0.0₃0693'
$ws.Range("E20").Value = '  -2.35%  '
$ws.Range("D21").Formula = '="7.32"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -2.68%  '
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Formula = '="3.86"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  -6.58%  '
$ws.Range("D24").Formula = '="9.14"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  -3.00%  '
$ws.Range("D25").Formula = '="2.06"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +5.93%  '
$ws.Range("D26").Formula = '="150.91"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  -0.65%  '
$ws.Range("D27").Formula = '="14.95"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  -2.52%  '
$ws.Range("E28").Value = '  -3.02%  '
$ws.Range("E29").Value = '  -3.40%  '
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("E31").Value = '  -1.85%  '
$ws.Range("E32").Value = '  -3.60%  '
$ws.Range("D33").Formula = '="3.19"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  -1.63%  '
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("D35").Value = '1.391.90'
$ws.Range("E35").Value = '  -1.94%  '
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("E37").Value = '  -4.01%  '
$ws.Range("E38").Value = '  +1.28%  '
$ws.Range("D39").Formula = '="2.57"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("D40").Formula = '="0.0165"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -1.18%  '
$ws.Range("D41").Formula = '="0.535"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -1.57%  '
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").Formula = '="0.786"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  -3.40%  '
$ws.Range("D44").Formula = '="1.86"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +1.06%  '
$ws.Range("D45").Formula = '="5.50"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  -4.38%  '
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").Formula = '="62.60"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  -3.28%  '
$ws.Range("D48").Value = '1.702.85'
$ws.Range("E48").Value = '  -1.50%  '
$ws.Range("D49").Formula = '="86.27"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  -1.89%  '
$ws.Range("E50").Value = '  -1.21%  '
$ws.Range("E51").Value = '  -0.85%  '

$excel.CutCopyMode = 0

